# Data updates July release part 1
#
# Updates the APT_ATC_PRE_APT workbook:
#  - Bumps the release date in the metadata header.
#  - Fills in previously-blank ATC pre-departure + all-causes figures for
#    Munich (EDDM, row 12) and Stuttgart (EDDS, row 13), extending the
#    existing shared "Tot. dep. delay" formula range down through row 17.
#  - Restates Toulouse-Blagnac (row 28) and Bucharest/Otopeni (row 45)
#    after additional data came in.
#  - Logs both changes on the "Change Log" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APT_ATC_PRE_APT")
$log = $wb.Worksheets.Item("Change Log")

# --- Metadata: release date (B2) -------------------------------------------
$ws.Range("B2").Value = 45477

# --- Munich (EDDM), row 12: add # departures / pre-departure delay ---------
$ws.Range("D12").Value = 149830
$ws.Range("E12").Value = 735

# --- Stuttgart (EDDS), row 13: add all figures ------------------------------
$ws.Range("D13").Value = 42730
$ws.Range("E13").Value = 2558
$ws.Range("G13").Value = 42730
$ws.Range("H13").Value = 613952.62

# Recreate the "ATC dep. delay (min./dep.)" shared formula across I6:I17 now
# that row 13 has real G/H data (previously split into I6:I12 / I14:I17).
$ws.Range("I6:I17").Formula = "=H6/G6"

# Recreate the "dep. delay (min./dep.)" shared formula across F11:F14 now
# that rows 12 and 13 have real D/E data.
$ws.Range("F11:F14").Formula = "=E11/D11"

# --- Toulouse-Blagnac (LFBO), row 28: updated figures -----------------------
$ws.Range("D28").Value = 34229
$ws.Range("E28").Value = 11475

# --- Bucharest/Otopeni (LROP), row 45: updated figure -----------------------
$ws.Range("D45").Value = 51003

# --- Change Log sheet --------------------------------------------------------
$log.Columns.Item(2).ColumnWidth = 8.88

$log.Range("A2").Value = 45477
$log.Range("B2").Value = "EDDM, EDDS"
$log.Range("C2").Value = 2023
$log.Range("D2").Value = "ATC pre-departure delays added for EDDM and EDDS, All pre-departure delays added for EDDS"

# Row 3 reuses row 2's styling (same column formats), then gets its own data.
$log.Range("A2:D2").Copy()
$log.Range("A3:D3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$log.Rows.Item(3).RowHeight = 14.25

$log.Range("A3").Value = 45477
$log.Range("B3").Value = "LFBO"
$log.Range("C3").Value = 2023
$log.Range("D3").Value = "Figures updated with additional data (old 31,494 ; 10,577)`n"

$excel.CutCopyMode = $false
